$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 337, shifting existing rows 337:384 down to 339:386
$ws.Rows("337:338").Insert()

# Populate new row 337 (Choclero, Primera, Region de O'Higgins)
$ws.Range("A337").Value = 4
$ws.Range("B337").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C337").Value = "Los Lagos"
$ws.Range("D337").Value = 45034
$ws.Range("E337").Value = 10
$ws.Range("F337").Value = 100112024
$ws.Range("G337").Value = "Choclo"
$ws.Range("H337").Value = "Choclero"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 3000
$ws.Range("K337").Value = 500
$ws.Range("L337").Value = 500
$ws.Range("M337").Value = 500
$ws.Range("N337").Value = "`$/unidad"
$ws.Range("O337").Value = "Región de O'Higgins"
$ws.Range("P337").Value = 500
$ws.Range("Q337").Value = 1
$ws.Range("R337").Value = "Hortaliza"

# Populate new row 338 (Choclero, Segunda, Region de O'Higgins)
$ws.Range("A338").Value = 4
$ws.Range("B338").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C338").Value = "Los Lagos"
$ws.Range("D338").Value = 45034
$ws.Range("E338").Value = 10
$ws.Range("F338").Value = 100112024
$ws.Range("G338").Value = "Choclo"
$ws.Range("H338").Value = "Choclero"
$ws.Range("I338").Value = "Segunda"
$ws.Range("J338").Value = 3000
$ws.Range("K338").Value = 480
$ws.Range("L338").Value = 480
$ws.Range("M338").Value = 480
$ws.Range("N338").Value = "`$/unidad"
$ws.Range("O338").Value = "Región de O'Higgins"
$ws.Range("P338").Value = 480
$ws.Range("Q338").Value = 1
$ws.Range("R338").Value = "Hortaliza"

# Ensure the date cells keep the existing date number format (style index 2 used throughout column D)
$ws.Range("D337").NumberFormat = $ws.Range("D339").NumberFormat
$ws.Range("D338").NumberFormat = $ws.Range("D339").NumberFormat
